$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> @{ D = "newPrice"; E = "newVolume" }
# Only entries present in a given row's hashtable are updated, matching
# the diff exactly (many rows only touch E, some only D+E).
$updates = @{
    2  = @{ D = "66.312.81";  E = "  -0.90%  " }
    3  = @{ D = "3.534.51";   E = "  +0.62%  " }
    4  = @{                   E = "  +0.09%  " }
    5  = @{ D = "607.32";     E = "  -0.22%  " }
    6  = @{ D = "143.94";     E = "  -2.57%  " }
    7  = @{ D = "3.533.37";   E = "  +0.59%  " }
    8  = @{                   E = "  +0.06%  " }
    9  = @{                   E = "  +0.41%  " }
    10 = @{ D = "8.10";       E = "  +1.12%  " }
    11 = @{                   E = "  -4.20%  " }
    12 = @{                   E = "  -3.00%  " }
    13 = @{ D = "4.134.08";   E = "  +0.69%  " }
    14 = @{                   E = "  -4.57%  " }
    15 = @{ D = "30.24";      E = "  -5.52%  " }
    16 = @{ D = "3.533.89";   E = "  +0.75%  " }
    17 = @{ D = "66.379.55";  E = "  -0.82%  " }
    18 = @{                   E = "  -0.69%  " }
    19 = @{                   E = "  +2.06%  " }
    20 = @{                   E = "  -3.84%  " }
    21 = @{ D = "14.92";      E = "  -2.87%  " }
    22 = @{                   E = "  -2.86%  " }
    23 = @{                   E = "  -1.42%  " }
    24 = @{ D = "78.78";      E = "  -1.11%  " }
    25 = @{ D = "3.676.44";   E = "  +0.79%  " }
    27 = @{                   E = "  -0.26%  " }
    28 = @{                   E = "  -1.66%  " }
    29 = @{ D = "9.22";       E = "  -5.74%  " }
    30 = @{                   E = "  -1.74%  " }
    31 = @{                   E = "  -0.11%  " }
    32 = @{ D = "1.49";       E = "  -7.33%  " }
    33 = @{                   E = "  -4.60%  " }
    34 = @{ D = "25.27";      E = "  -1.19%  " }
    35 = @{ D = "3.525.88" }
    37 = @{                   E = "  -3.21%  " }
    38 = @{ D = "7.81";       E = "  -3.26%  " }
    39 = @{                   E = "  -5.63%  " }
    40 = @{ D = "1.00";       E = "  +0.10%  " }
    41 = @{ D = "173.07";     E = "  -1.20%  " }
    42 = @{                   E = "  -4.33%  " }
    43 = @{ D = "5.17";       E = "  -4.25%  " }
    44 = @{ D = "0.892";      E = "  -0.31%  " }
    45 = @{                   E = "  -7.81%  " }
    46 = @{ D = "45.46";      E = "  -1.39%  " }
    47 = @{ D = "1.23";       E = "  -2.37%  " }
    48 = @{ D = "26.06";      E = "  -7.43%  " }
    49 = @{ D = "2.40";       E = "  -1.86%  " }
    50 = @{ D = "7.13";       E = "  -4.52%  " }
    51 = @{                   E = "  -4.90%  " }
}

# These "Price" / "Volume" columns are stored as plain text in the
# workbook even when the text looks like a number (e.g. "607.32").
# Assigning a bare numeric-looking string via .Value lets Excel's
# automatic type detection turn it into a real number, which would
# change the cell's stored type. Prefixing with a leading apostrophe
# forces Excel to keep/treat the entry as literal text (exactly like a
# user typing '607.32 into a cell) while still storing the clean
# "607.32" text without the apostrophe itself. Values such as
# "66.312.81" (two dots) are never valid numbers, so they are always
# safe to assign directly.

foreach ($rowNum in $updates.Keys) {
    $rowData = $updates[$rowNum]

    if ($rowData.ContainsKey("D")) {
        $dText = $rowData["D"]
        if ($dText -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
            $ws.Range("D$rowNum").Value = "'" + $dText
        } else {
            $ws.Range("D$rowNum").Value = $dText
        }
    }

    if ($rowData.ContainsKey("E")) {
        $eText = $rowData["E"]
        if ($eText -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
            $ws.Range("E$rowNum").Value = "'" + $eText
        } else {
            $ws.Range("E$rowNum").Value = $eText
        }
    }
}
